$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.799.51'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.08%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.291.89'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.22%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.37%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '113.81'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +15.50%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '268.70'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.76%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.21%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.40%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '48.84'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +7.70%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0942'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.70%  '

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +14.42%  '

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.49%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.84'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.24%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.634.40'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.27%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.871'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.67%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.293.07'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.04%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.696.63'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.18%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.74%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.00'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +12.64%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.26'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.09%  '

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.99%  '

$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '233.14'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.08%  '

$ws.Range("B24").Value = 'InternetComputer(DFINITY)'
$ws.Range("C24").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.83'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +7.53%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.77%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.74'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +3.84%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.06%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '43.49'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +14.04%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.52%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.05%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.32%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '173.40'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.85%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0935'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +4.33%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '21.62'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.13%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.68'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +4.64%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.67%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.80'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.17%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.58%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.106'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.59%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.80'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +7.61%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '14.70'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +20.63%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '74.53'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +14.96%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.40'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.76%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.242'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +2.20%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.32'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +21.10%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.17%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.41'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.73%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.59%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '103.04'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +4.56%  '

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +3.87%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.27%  '
